$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert a new row above the current row 18 (AuthorisationStatuses), pushing
# AuthorisationStatuses / ApprovalType / SelectAllTransactions / FromDate /
# ToDate down by one row each.
$ws.Range("A18").EntireRow.Insert()

# Restore the standard row height/format on the freshly inserted row so it
# matches the rest of the sheet's formatted rows.
$ws.Rows.Item(18).RowHeight = 14.25

# Populate the newly inserted row 18 with the new "MaxPageRuns" setting.
$ws.Range("A18").Value = "MaxPageRuns"
$ws.Range("B18").Value = 3

# Reflect the UI state recorded after the edit: selection moved to C23
# (the now-shifted "ToDate" row).
$ws.Range("C23").Select() | Out-Null
